# Update the NATMI LR-pairs sheet (Ccl5-Cxcr3) with recomputed TPM-based
# statistics. Columns A-D (Sending cluster, Ligand symbol, Receptor symbol,
# Target cluster) are unchanged; columns E..T (all the numeric expression /
# specificity / edge-weight metrics) are refreshed with the new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row number -> new values for columns E..T (in column order).
$newValues = @{
  2  = @(3, 1, 1.141514, 3.424542, 0.1286876091642341, 0.1286876091642341, 1, 0.3333333333333333, 0.065175, 0.195525, 0.009404016458916581, 0.009404016458916581, 0.07439817495000001, 0.6695835745500001, 0.001210180394639082, 0.001210180394639082)
  3  = @(3, 1, 1.141514, 3.424542, 0.1286876091642341, 0.1286876091642341, 3, 1, 6.718514333333332, 20.155543, 0.969405744075698, 0.969405744075698, 7.669278170700666, 69.02350353630599, 0.124750507515177, 0.124750507515177)
  4  = @(3, 1, 1.141514, 3.424542, 0.1286876091642341, 0.1286876091642341, 1, 0.3333333333333333, 0.14686, 0.44058, 0.02119023946538534, 0.02119023946538533, 0.16764274604, 1.50878471436, 0.002726921254418038, 0.002726921254418037)
  5  = @(3, 1, 1.348859666666667, 4.046578999999999, 0.1520625464088912, 0.1520625464088912, 1, 0.3333333333333333, 0.065175, 0.195525, 0.009404016458916581, 0.009404016458916581, 0.087911928775, 0.7912073589749999, 0.00142999868921398, 0.00142999868921398)
  6  = @(3, 1, 1.348859666666667, 4.046578999999999, 0.1520625464088912, 0.1520625464088912, 3, 1, 6.718514333333332, 20.155543, 0.969405744075698, 0.969405744075698, 9.062333004155221, 81.56099703739699, 0.1474103059475566, 0.1474103059475566)
  7  = @(3, 1, 1.348859666666667, 4.046578999999999, 0.1520625464088912, 0.1520625464088912, 1, 0.3333333333333333, 0.14686, 0.44058, 0.02119023946538534, 0.02119023946538533, 0.1980935306466667, 1.78284177582, 0.003222241772120677, 0.003222241772120676)
  8  = @(3, 1, 2.436302666666667, 7.308908000000001, 0.2746545074118945, 0.2746545074118945, 1, 0.3333333333333333, 0.065175, 0.195525, 0.009404016458916581, 0.009404016458916581, 0.1587860263, 1.4290742367, 0.002582855508217082, 0.002582855508217082)
  9  = @(3, 1, 2.436302666666667, 7.308908000000001, 0.2746545074118945, 0.2746545074118945, 3, 1, 6.718514333333332, 20.155543, 0.969405744075698, 0.969405744075698, 16.36833438633822, 147.315009477044, 0.2662516571213719, 0.2662516571213719)
  10 = @(3, 1, 2.436302666666667, 7.308908000000001, 0.2746545074118945, 0.2746545074118945, 1, 0.3333333333333333, 0.14686, 0.44058, 0.02119023946538534, 0.02119023946538533, 0.3577954096266667, 3.22015868664, 0.005819994782305497, 0.005819994782305496)
  11 = @(1, 0.3333333333333333, 0.4110173333333333, 1.233052, 0.04633568922652349, 0.04633568922652349, 1, 0.3333333333333333, 0.065175, 0.195525, 0.009404016458916581, 0.009404016458916581, 0.0267880547, 0.2410924923, 0.0004357415841214706, 0.0004357415841214706)
  12 = @(1, 0.3333333333333333, 0.4110173333333333, 1.233052, 0.04633568922652349, 0.04633568922652349, 3, 1, 6.718514333333332, 20.155543, 0.969405744075698, 0.969405744075698, 2.761425845248444, 24.852832607236, 0.0449180832918983, 0.0449180832918983)
  13 = @(1, 0.3333333333333333, 0.4110173333333333, 1.233052, 0.04633568922652349, 0.04633568922652349, 1, 0.3333333333333333, 0.14686, 0.44058, 0.02119023946538534, 0.02119023946538533, 0.06036200557333334, 0.5432580501600001, 0.0009818643505037081, 0.0009818643505037081)
  14 = @(3, 1, 2.206596333333333, 6.619789, 0.2487587594433638, 0.2487587594433639, 1, 0.3333333333333333, 0.065175, 0.195525, 0.009404016458916581, 0.009404016458916581, 0.143814916025, 1.294334244225, 0.002339331468105064, 0.002339331468105064)
  15 = @(3, 1, 2.206596333333333, 6.619789, 0.2487587594433638, 0.2487587594433639, 3, 1, 6.718514333333332, 20.155543, 0.969405744075698, 0.969405744075698, 14.82504909338078, 133.425441840427, 0.2411481702935417, 0.2411481702935417)
  16 = @(3, 1, 2.206596333333333, 6.619789, 0.2487587594433638, 0.2487587594433639, 1, 0.3333333333333333, 0.14686, 0.44058, 0.02119023946538534, 0.02119023946538533, 0.3240607375133334, 2.91654663762, 0.005271257681717066, 0.005271257681717065)
  17 = @(3, 1, 1.326136666666667, 3.97841, 0.1495008883450928, 0.1495008883450928, 1, 0.3333333333333333, 0.065175, 0.195525, 0.009404016458916581, 0.009404016458916581, 0.08643095725, 0.77787861525, 0.001405908814619902, 0.001405908814619902)
  18 = @(3, 1, 1.326136666666667, 3.97841, 0.1495008883450928, 0.1495008883450928, 3, 1, 6.718514333333332, 20.155543, 0.969405744075698, 0.969405744075698, 8.909668202958889, 80.18701382662999, 0.1449270199061525, 0.1449270199061525)
  19 = @(3, 1, 1.326136666666667, 3.97841, 0.1495008883450928, 0.1495008883450928, 1, 0.3333333333333333, 0.14686, 0.44058, 0.02119023946538534, 0.02119023946538533, 0.1947564308666667, 1.7528078778, 0.003167959624320351, 0.003167959624320351)
}

# Column E is worksheet column index 5; write each row's values into E..T.
foreach ($rowNum in $newValues.Keys) {
  $rowVals = $newValues[$rowNum]
  for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $colIndex = 5 + $i
    $ws.Cells.Item([int]$rowNum, $colIndex).Value2 = $rowVals[$i]
  }
}
